# Append the standard "new blank Word document" placeholder paragraphs
# after the existing "Test" paragraph, and push the (hidden) _GoBack
# bookmark down into its own empty paragraph at the very end of the
# document - mirroring what happens when a user places the cursor at
# the end of the document and types several new paragraphs.

$d = $word.ActiveDocument

# The _GoBack bookmark is hidden from the Bookmarks collection/iteration
# but can still be addressed by name. Remove it first so the "Test"
# paragraph can be left holding only the "Test" run, then we re-create
# the bookmark collapsed at the new document end once all the new
# paragraphs/text are in place.
$hadGoBack = $false
try {
    $goBack = $d.Bookmarks("_GoBack")
    $hadGoBack = $true
    $goBack.Delete()
} catch {
    $hadGoBack = $false
}

$paragraphTexts = @(
    "Video bietet eine leistungsstarke Möglichkeit zur Unterstützung Ihres Standpunkts. Wenn Sie auf `"Onlinevideo`" klicken, können Sie den Einbettungscode für das Video einfügen, das hinzugefügt werden soll. Sie können auch ein Stichwort eingeben, um online nach dem Videoclip zu suchen, der optimal zu Ihrem Dokument passt.",
    "Damit Ihr Dokument ein professionelles Aussehen erhält, stellt Word einander ergänzende Designs für Kopfzeile, Fußzeile, Deckblatt und Textfelder zur Verfügung. Beispielsweise können Sie ein passendes Deckblatt mit Kopfzeile und Randleiste hinzufügen. Klicken Sie auf `"Einfügen`", und wählen Sie dann die gewünschten Elemente aus den verschiedenen Katalogen aus.",
    "Designs und Formatvorlagen helfen auch dabei, die Elemente Ihres Dokuments aufeinander abzustimmen. Wenn Sie auf `"Design`" klicken und ein neues Design auswählen, ändern sich die Grafiken, Diagramme und SmartArt-Grafiken so, dass sie dem neuen Design entsprechen. Wenn Sie Formatvorlagen anwenden, ändern sich die Überschriften passend zum neuen Design.",
    "Sparen Sie Zeit in Word dank neuer Schaltflächen, die angezeigt werden, wo Sie sie benötigen. Zum Ändern der Weise, in der sich ein Bild in Ihr Dokument einfügt, klicken Sie auf das Bild. Dann wird eine Schaltfläche für Layoutoptionen neben dem Bild angezeigt Beim Arbeiten an einer Tabelle klicken Sie an die Position, an der Sie eine Zeile oder Spalte hinzufügen möchten, und klicken Sie dann auf das Pluszeichen.",
    "Auch das Lesen ist bequemer in der neuen Leseansicht. Sie können Teile des Dokuments reduzieren und sich auf den gewünschten Text konzentrieren. Wenn Sie vor dem Ende zu lesen aufhören müssen, merkt sich Word die Stelle, bis zu der Sie gelangt sind – sogar auf einem anderen Gerät."
)

# Paragraph layout being built after paragraph 1 ("Test"):
#   2           -> blank paragraph
#   3 .. 7      -> the five text paragraphs above
#   8 (last)    -> blank paragraph that will hold the bookmark
$totalNewParagraphs = 1 + $paragraphTexts.Count + 1
for ($i = 0; $i -lt $totalNewParagraphs; $i++) {
    $endRange = $d.Content
    $endRange.Collapse(0)
    $endRange.InsertParagraphAfter()
}

for ($i = 0; $i -lt $paragraphTexts.Count; $i++) {
    $textParagraph = $d.Paragraphs($i + 3)
    $textParagraph.Range.Text = $paragraphTexts[$i]
}

if ($hadGoBack) {
    # Re-create _GoBack collapsed at the new document end. Adding a
    # bookmark directly on a zero-length range at a paragraph boundary
    # is unreliable, so bracket a throwaway placeholder character
    # instead, bookmark that one-character range, then delete just the
    # placeholder - leaving the bookmark correctly collapsed in place.
    $endRange = $d.Content
    $endRange.Collapse(0)
    $boundary = $endRange.Start
    $endRange.InsertAfter("X")

    $placeholderRange = $d.Range($boundary - 1, $boundary)
    $d.Bookmarks.Add("_GoBack", $placeholderRange)

    $placeholderRange = $d.Range($boundary - 1, $boundary)
    $placeholderRange.Delete()
}
